# Auto-generated edit script: applies the Masamune_Profits.xlsx leve-price refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (matching the scheduled-runner diff).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 19465.303
$ws.Range("I132").Value = 2605.6135
$ws.Range("K132").Value = 7816.8405
$ws.Range("M132").Value = -5286.8405

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 52984.5
$ws.Range("J44").Value = 52984.5
$ws.Range("L44").Value = 52984.5
$ws.Range("N44").Value = -53960.5
$ws.Range("H61").Value = 2121.6584
$ws.Range("I61").Value = 1572.1538
$ws.Range("K61").Value = 1572.1538
$ws.Range("M61").Value = -1360.1538
$ws.Range("H74").Value = 2137.9333
$ws.Range("I74").Value = 1782.1428
$ws.Range("J74").Value = 2968.111
$ws.Range("K74").Value = 1782.1428
$ws.Range("L74").Value = 2968.111
$ws.Range("M74").Value = -908.1428000000001
$ws.Range("N74").Value = -4716.111
$ws.Range("H77").Value = 2137.9333
$ws.Range("I77").Value = 1782.1428
$ws.Range("J77").Value = 2968.111
$ws.Range("K77").Value = 8910.714
$ws.Range("L77").Value = 14840.555
$ws.Range("M77").Value = -4542.714
$ws.Range("N77").Value = -23576.555
$ws.Range("H132").Value = 13516425
$ws.Range("I132").Value = 26317796
$ws.Range("J132").Value = 3866.889
$ws.Range("K132").Value = 78953388
$ws.Range("L132").Value = 11600.667
$ws.Range("M132").Value = -78950858
$ws.Range("N132").Value = -16660.667
$ws.Range("H136").Value = 2121.6584
$ws.Range("I136").Value = 1572.1538
$ws.Range("K136").Value = 4716.4614
$ws.Range("M136").Value = -2166.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2206
$ws.Range("I134").Value = 1720.5625
$ws.Range("K134").Value = 5161.6875
$ws.Range("M134").Value = -2626.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5468731.5
$ws.Range("I31").Value = 2267
$ws.Range("J31").Value = 7580774.5
$ws.Range("K31").Value = 2267
$ws.Range("L31").Value = 7580774.5
$ws.Range("M31").Value = -1972
$ws.Range("N31").Value = -7581364.5
$ws.Range("H34").Value = 5468731.5
$ws.Range("I34").Value = 2267
$ws.Range("J34").Value = 7580774.5
$ws.Range("K34").Value = 2267
$ws.Range("L34").Value = 7580774.5
$ws.Range("M34").Value = -2065
$ws.Range("N34").Value = -7581178.5
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H105").Value = 2405.75
$ws.Range("I105").Value = 2675.9375
$ws.Range("J105").Value = 1325
$ws.Range("K105").Value = 2675.9375
$ws.Range("L105").Value = 1325
$ws.Range("M105").Value = -928.9375
$ws.Range("N105").Value = -4819
$ws.Range("H132").Value = 42399.027
$ws.Range("I132").Value = 1960.125
$ws.Range("J132").Value = 130629.37
$ws.Range("K132").Value = 5880.375
$ws.Range("L132").Value = 391888.11
$ws.Range("M132").Value = -3350.375
$ws.Range("N132").Value = -396948.11
$ws.Range("H134").Value = 379601.34
$ws.Range("I134").Value = 1036.0312
$ws.Range("K134").Value = 3108.0936
$ws.Range("M134").Value = -573.0935999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 7058.5835
$ws.Range("I34").Value = 425.75
$ws.Range("J34").Value = 10375
$ws.Range("K34").Value = 1277.25
$ws.Range("L34").Value = 31125
$ws.Range("M34").Value = -1193.25
$ws.Range("N34").Value = -31293
$ws.Range("H55").Value = 3136.3635
$ws.Range("J55").Value = 3136.3635
$ws.Range("L55").Value = 9409.0905
$ws.Range("N55").Value = -9763.0905
$ws.Range("H17").Value = 1050
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 1160
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 3480
$ws.Range("M17").Value = -1331
$ws.Range("N17").Value = -3818
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H63").Value = 4220
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 4220
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 12660
$ws.Range("N63").Value = -14158
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 4220
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 4220
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 37980
$ws.Range("N66").Value = -45468
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 37895
$ws.Range("J52").Value = 37895
$ws.Range("L52").Value = 37895
$ws.Range("N52").Value = -38413
$ws.Range("H102").Value = 1289.8695
$ws.Range("I102").Value = 1124.4667
$ws.Range("J102").Value = 1600
$ws.Range("K102").Value = 1124.4667
$ws.Range("L102").Value = 1600
$ws.Range("M102").Value = 497.5333000000001
$ws.Range("N102").Value = -4844
$ws.Range("H132").Value = 2746.7778
$ws.Range("I132").Value = 1995.2727
$ws.Range("K132").Value = 5985.8181
$ws.Range("M132").Value = -3455.8181
$ws.Range("H135").Value = 67824.22
$ws.Range("J135").Value = 67824.22
$ws.Range("L135").Value = 67824.22
$ws.Range("N135").Value = -77964.22

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1239
$ws.Range("I93").Value = 1393.5
$ws.Range("J93").Value = 1215.2307
$ws.Range("K93").Value = 1393.5
$ws.Range("L93").Value = 1215.2307
$ws.Range("M93").Value = -145.5
$ws.Range("N93").Value = -3711.2307
$ws.Range("H132").Value = 3842.5483
$ws.Range("I132").Value = 3077.875
$ws.Range("J132").Value = 4658.2
$ws.Range("K132").Value = 9233.625
$ws.Range("L132").Value = 13974.6
$ws.Range("M132").Value = -6703.625
$ws.Range("N132").Value = -19034.6
$ws.Range("H136").Value = 1776.7428
$ws.Range("I136").Value = 1368.037
$ws.Range("J136").Value = 3156.125
$ws.Range("K136").Value = 4104.111
$ws.Range("L136").Value = 9468.375
$ws.Range("M136").Value = -1554.111
$ws.Range("N136").Value = -14568.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2094.8
$ws.Range("I132").Value = 1523.2667
$ws.Range("J132").Value = 2666.3333
$ws.Range("K132").Value = 4569.800099999999
$ws.Range("L132").Value = 7998.999899999999
$ws.Range("M132").Value = -2039.800099999999
$ws.Range("N132").Value = -13058.9999
$ws.Range("H136").Value = 182700.56
$ws.Range("I136").Value = 238779.23
$ws.Range("J136").Value = 1523.3077
$ws.Range("K136").Value = 716337.6900000001
$ws.Range("L136").Value = 4569.9231
$ws.Range("M136").Value = -713787.6900000001
$ws.Range("N136").Value = -9669.9231
